$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'59.192.08"
$ws.Range('E2').Value = '  +1.51%  '
$ws.Range('D3').Value = "'2.591.78"
$ws.Range('E3').Value = '  -0.23%  '
$ws.Range('D5').Value = "'524.55"
$ws.Range('E5').Value = '  +0.49%  '
$ws.Range('D6').Value = "'139.97"
$ws.Range('E6').Value = '  -2.73%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  -0.81%  '
$ws.Range('D9').Value = "'2.603.58"
$ws.Range('E9').Value = '  -0.58%  '
$ws.Range('E10').Value = '  +0.39%  '
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('D12').Value = "'0.332"
$ws.Range('E12').Value = '  -2.01%  '
$ws.Range('E13').Value = '  +2.97%  '
$ws.Range('D14').Value = "'3.049.31"
$ws.Range('E14').Value = '  -0.20%  '
$ws.Range('D15').Value = "'59.105.64"
$ws.Range('E15').Value = '  +1.47%  '
$ws.Range('D16').Value = "'20.54"
$ws.Range('E16').Value = '  +0.52%  '
$ws.Range('B17').Value = "'ShibaInu"
$ws.Range('C17').Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range('D17').Value = "'0.0000133"
$ws.Range('E17').Value = '  -0.39%  '
$ws.Range('B18').Value = "'WrappedEther"
$ws.Range('C18').Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range('D18').Value = "'2.575.59"
$ws.Range('E18').Value = '  -1.64%  '
$ws.Range('D19').Value = "'342.17"
$ws.Range('E19').Value = '  +0.62%  '
$ws.Range('E20').Value = '  -0.83%  '
$ws.Range('D21').Value = "'10.13"
$ws.Range('E21').Value = '  -1.40%  '
$ws.Range('E22').Value = '  +0.37%  '
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').Value = "'66.66"
$ws.Range('E25').Value = '  +0.47%  '
$ws.Range('E26').Value = '  +0.53%  '
$ws.Range('D27').Value = "'0.998"
$ws.Range('E27').Value = '  -0.17%  '
$ws.Range('D28').Value = "'7.09"
$ws.Range('E28').Value = '  +0.92%  '
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').Value = "'0.0₃0727"
$ws.Range('E30').Value = '  -3.03%  '
$ws.Range('D31').Value = "'5.94"
$ws.Range('E31').Value = '  -4.29%  '
$ws.Range('E32').Value = '  +0.36%  '
$ws.Range('D33').Value = "'18.76"
$ws.Range('E33').Value = '  -0.22%  '
$ws.Range('D34').Value = "'149.20"
$ws.Range('E34').Value = '  -0.50%  '
$ws.Range('D35').Value = "'3.99"
$ws.Range('E35').Value = '  -0.75%  '
$ws.Range('E36').Value = '  -1.08%  '
$ws.Range('D37').Value = "'36.80"
$ws.Range('E37').Value = '  +2.27%  '
$ws.Range('E38').Value = '  +1.45%  '
$ws.Range('D39').Value = "'0.832"
$ws.Range('E39').Value = '  -4.66%  '
$ws.Range('E40').Value = '  -6.86%  '
$ws.Range('D41').Value = "'3.53"
$ws.Range('E41').Value = '  -0.22%  '
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('D43').Value = "'272.64"
$ws.Range('E43').Value = '  +0.50%  '
$ws.Range('E44').Value = '  -0.73%  '
$ws.Range('E45').Value = '  +0.93%  '
$ws.Range('D46').Value = "'0.0953"
$ws.Range('E46').Value = '  -0.69%  '
$ws.Range('E47').Value = '  -1.61%  '
$ws.Range('D48').Value = "'18.46"
$ws.Range('E48').Value = '  -1.81%  '
$ws.Range('D49').Value = "'1.969.82"
$ws.Range('E49').Value = '  -0.17%  '
$ws.Range('B50').Value = "'RenderToken"
$ws.Range('C50').Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range('D50').Value = "'4.64"
$ws.Range('E50').Value = '  +0.50%  '
$ws.Range('B51').Value = "'VeChain"
$ws.Range('C51').Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range('D51').Value = "'0.0223"
$ws.Range('E51').Value = '  +0.32%  '
